# Update rules in DiscountRules.xlsx
# - Row 25, column A is cleared (blank)
# - Row 25, column C gets a new test marker string
# - Rows 26 and 27 (the old "3rd condition"/"Rule 3"/"3rd Rule" rows) are removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(25, 1).Value = ""
$ws.Cells.Item(25, 3).Value = "Test - 09302025 1056"

$ws.Rows("26:27").Delete()
